# Daily attendance processing - 2026-01-20 23:02:02
# Normalize the "Recorded By" (column G) entries so the recorder's email
# address is listed first, followed by "System", instead of the other
# way around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$col = $ws.Columns.Item(7)  # Column G = "Recorded By"

$first = $col.Find($oldValue)
if ($first -ne $null) {
    $firstAddress = $first.Address()
    $current = $first
    do {
        $current.Value2 = $newValue
        $current = $col.FindNext($current)
    } while ($current -ne $null -and $current.Address() -ne $firstAddress)
}
